$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# StatOutput sheet: add the stat-query result table (header row + value row)
# ---------------------------------------------------------------------------
$statOutput = $wb.Worksheets.Item("StatOutput")

$statOutput.Range("A1").Value = "number_of_files"
$statOutput.Range("B1").Value = "number_of_sample"
$statOutput.Range("C1").Value = "number_of_cases"
$statOutput.Range("D1").Value = "number_of_study"

# The counts look like numbers, so force them to be stored as text (shared
# strings) the same way the source workbook does, using the apostrophe
# text-prefix so Excel doesn't coerce them into numeric cells.
$statOutput.Range("A2").Value = "'0"
$statOutput.Range("B2").Value = "'0"
$statOutput.Range("C2").Value = "'18"
$statOutput.Range("D2").Value = "'1"

# ---------------------------------------------------------------------------
# StatOutput_Message sheet: the "empty cypher" error run is replaced by a
# second successful connection-info block that now carries the real stats
# cypher query text instead of an empty string.
# ---------------------------------------------------------------------------
$statMsg = $wb.Worksheets.Item("StatOutput_Message")

# Row 11 used to hold "Cypher query should not be an empty string" - that
# run no longer errors, so the row disappears and everything below shifts
# up by one.
$statMsg.Rows(11).Delete()

$statsCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN['Lymphoma']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# After the row shift, row 18 (which held the empty cypher placeholder) now
# gets the real stats cypher text.
$statMsg.Range("A18").Value = $statsCypher

# Row 21 (the trailing output filepath line) is restored after the shift.
$statMsg.Range("A21").Value = "C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC05_Canine_Filter_Diagnosis-Lymphoma_Neo4jData.xlsx"
